$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from the row above so the new row matches existing styling exactly
$ws.Range("A28:D28").Copy()
$ws.Range("A29:D29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new row 29 data: 485. Max Consecutive Ones (Java, dated 2023-03-18 / serial 45003)
$ws.Range("A29").Value = 485
$ws.Range("B29").Value = "Max Consecutive Ones"
$ws.Range("C29").Value = "Java "
$ws.Range("D29").Value = 45003

# Update the visible selection to match the saved state
$ws.Range("K29").Select()
